# Update countries & provincias Spain
# Updates country case numbers and re-orders a few country rows
# (Mauricio, Trinidad y Tobago, Mali moved up one position each)
# plus refreshes the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in row 1
$ws.Range("A1").Value2 = "Datos actualizados a 7 de Abril de 2020 a las 16:22"

$ws.Range("A4").Value2 = 'Estados Unidos'
$ws.Range("B4").Value2 = 369179
$ws.Range("C4").Value2 = 2175
$ws.Range("D4").Value2 = 19843
$ws.Range("E4").Value2 = 338323
$ws.Range("F4").Value2 = 8983
$ws.Range("G4").Value2 = 142
$ws.Range("H4").Value2 = 11013
$ws.Range("A13").Value2 = 'Suiza'
$ws.Range("B13").Value2 = 22242
$ws.Range("C13").Value2 = 585
$ws.Range("D13").Value2 = 8056
$ws.Range("E13").Value2 = 13375
$ws.Range("F13").Value2 = 391
$ws.Range("G13").Value2 = 46
$ws.Range("H13").Value2 = 811
$ws.Range("A17").Value2 = 'Austria'
$ws.Range("B17").Value2 = 12547
$ws.Range("C17").Value2 = 250
$ws.Range("D17").Value2 = 4046
$ws.Range("E17").Value2 = 8258
$ws.Range("F17").Value2 = 243
$ws.Range("G17").Value2 = 23
$ws.Range("H17").Value2 = 243
$ws.Range("A51").Value2 = 'Sudafrica'
$ws.Range("B51").Value2 = 1749
$ws.Range("C51").Value2 = 63
$ws.Range("D51").Value2 = 95
$ws.Range("E51").Value2 = 1641
$ws.Range("F51").Value2 = 7
$ws.Range("G51").Value2 = 1
$ws.Range("H51").Value2 = 13
$ws.Range("A63").Value2 = 'Eslovenia'
$ws.Range("B63").Value2 = 1059
$ws.Range("C63").Value2 = 38
$ws.Range("D63").Value2 = 102
$ws.Range("E63").Value2 = 921
$ws.Range("F63").Value2 = 30
$ws.Range("G63").Value2 = 6
$ws.Range("H63").Value2 = 36
$ws.Range("A68").Value2 = 'Bielorrusia'
$ws.Range("B68").Value2 = 861
$ws.Range("C68").Value2 = 161
$ws.Range("D68").Value2 = 54
$ws.Range("E68").Value2 = 794
$ws.Range("F68").Value2 = 31
$ws.Range("G68").Value2 = 0
$ws.Range("H68").Value2 = 13
$ws.Range("A102").Value2 = 'Mauricio'
$ws.Range("B102").Value2 = 268
$ws.Range("C102").Value2 = 24
$ws.Range("D102").Value2 = 8
$ws.Range("E102").Value2 = 253
$ws.Range("F102").Value2 = 2
$ws.Range("G102").Value2 = 0
$ws.Range("H102").Value2 = 7
$ws.Range("A103").Value2 = 'Estado de Palestina'
$ws.Range("B103").Value2 = 260
$ws.Range("C103").Value2 = 6
$ws.Range("D103").Value2 = 24
$ws.Range("E103").Value2 = 235
$ws.Range("F103").Value2 = 0
$ws.Range("G103").Value2 = 0
$ws.Range("H103").Value2 = 1
$ws.Range("A104").Value2 = 'Niger'
$ws.Range("B104").Value2 = 253
$ws.Range("C104").Value2 = 0
$ws.Range("D104").Value2 = 26
$ws.Range("E104").Value2 = 217
$ws.Range("F104").Value2 = 0
$ws.Range("G104").Value2 = 0
$ws.Range("H104").Value2 = 10
$ws.Range("A105").Value2 = 'Vietnam'
$ws.Range("B105").Value2 = 249
$ws.Range("C105").Value2 = 4
$ws.Range("D105").Value2 = 123
$ws.Range("E105").Value2 = 126
$ws.Range("F105").Value2 = 8
$ws.Range("G105").Value2 = 0
$ws.Range("H105").Value2 = 0
$ws.Range("A110").Value2 = 'Georgia'
$ws.Range("B110").Value2 = 195
$ws.Range("C110").Value2 = 7
$ws.Range("D110").Value2 = 45
$ws.Range("E110").Value2 = 147
$ws.Range("F110").Value2 = 6
$ws.Range("G110").Value2 = 1
$ws.Range("H110").Value2 = 3
$ws.Range("A127").Value2 = 'Trinidad yTobago'
$ws.Range("B127").Value2 = 106
$ws.Range("C127").Value2 = 1
$ws.Range("D127").Value2 = 1
$ws.Range("E127").Value2 = 97
$ws.Range("F127").Value2 = 0
$ws.Range("G127").Value2 = 0
$ws.Range("H127").Value2 = 8
$ws.Range("A128").Value2 = 'Ruanda'
$ws.Range("B128").Value2 = 105
$ws.Range("C128").Value2 = 0
$ws.Range("D128").Value2 = 4
$ws.Range("E128").Value2 = 101
$ws.Range("F128").Value2 = 0
$ws.Range("G128").Value2 = 0
$ws.Range("H128").Value2 = 0
$ws.Range("A140").Value2 = 'Mali'
$ws.Range("B140").Value2 = 56
$ws.Range("C140").Value2 = 9
$ws.Range("D140").Value2 = 12
$ws.Range("E140").Value2 = 39
$ws.Range("F140").Value2 = 0
$ws.Range("G140").Value2 = 0
$ws.Range("H140").Value2 = 5
$ws.Range("A141").Value2 = 'Uganda'
$ws.Range("B141").Value2 = 52
$ws.Range("C141").Value2 = 0
$ws.Range("D141").Value2 = 0
$ws.Range("E141").Value2 = 52
$ws.Range("F141").Value2 = 0
$ws.Range("G141").Value2 = 0
$ws.Range("H141").Value2 = 0
$ws.Range("A142").Value2 = 'Etiopia'
$ws.Range("B142").Value2 = 52
$ws.Range("C142").Value2 = 8
$ws.Range("D142").Value2 = 4
$ws.Range("E142").Value2 = 46
$ws.Range("F142").Value2 = 1
$ws.Range("G142").Value2 = 0
$ws.Range("H142").Value2 = 2
